$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Split the run "||stay_on_the_page||" into "||stay_on_the_page" and "||"
#    (clean run split, no residual formatting) using the
#    insert-paragraph-mark-then-delete-it trick.
# ---------------------------------------------------------------------------
$r = $d.Content
$found = $r.Find.Execute("||stay_on_the_page||")
if (-not $found) {
    throw "Could not find '||stay_on_the_page||' in the document"
}
$matchEnd = $r.End
$splitPoint = $matchEnd - 2

$splitRange = $d.Range($splitPoint, $splitPoint)
$splitRange.InsertParagraphAfter()
$markRange = $d.Range($splitPoint, $splitPoint + 1)
$markRange.Delete()

# This paragraph (still holding "...||stay_on_the_page|| so that you can
# debug...provider.") keeps its original index in the Paragraphs collection
# because the insert+delete above nets out to zero new paragraphs.
$hostParaIndex = 28
$hostPara = $d.Paragraphs.Item($hostParaIndex)

# ---------------------------------------------------------------------------
# 2) Move the "_GoBack" bookmark out of this paragraph into its own, brand
#    new, otherwise-empty paragraph placed right after it.
# ---------------------------------------------------------------------------
$bm = $d.Bookmarks.Item("_GoBack")
$bm.Delete()

# Insert a brand new empty paragraph right after the host paragraph.
$hostPara.Range.InsertParagraphAfter()
$bmParaIndex = $hostParaIndex + 1
$bmPara = $d.Paragraphs.Item($bmParaIndex)

# Work around a zero-length-range-at-paragraph-boundary quirk in
# Bookmarks.Add by temporarily inserting a placeholder character, anchoring
# the bookmark next to it (no longer an ambiguous paragraph-boundary point),
# then removing the placeholder again.
$placeholder = $d.Range($bmPara.Range.Start, $bmPara.Range.Start)
$placeholder.InsertAfter("X")
$bmParaAfter = $d.Paragraphs.Item($bmParaIndex)
$bmAnchor = $d.Range($bmParaAfter.Range.Start, $bmParaAfter.Range.Start)
$d.Bookmarks.Add("_GoBack", $bmAnchor)
$xRange = $d.Range($bmParaAfter.Range.Start, $bmParaAfter.Range.Start + 1)
$xRange.Delete()

# ---------------------------------------------------------------------------
# 3) Add a brand new paragraph after the bookmark paragraph with the
#    "Here is some text ||side||||by_side|| Here is more text" content,
#    underlining the middle "||||" run.
# ---------------------------------------------------------------------------
$bmParaFinal = $d.Paragraphs.Item($bmParaIndex)
$bmParaFinal.Range.InsertParagraphAfter()
$newParaIndex = $bmParaIndex + 1
$newPara = $d.Paragraphs.Item($newParaIndex)

$fullText = "Here is some text ||side||||by_side|| Here is more text"
$insertNew = $d.Range($newPara.Range.Start, $newPara.Range.Start)
$insertNew.InsertAfter($fullText)

$base = $newPara.Range.Start

# Split "Here is some text ||" from "side" (plain/plain boundary -> use the
# insert-then-delete paragraph-mark trick again to avoid leftover rPr).
$b1 = $base + 20
$splitRange1 = $d.Range($b1, $b1)
$splitRange1.InsertParagraphAfter()
$markRange1 = $d.Range($b1, $b1 + 1)
$markRange1.Delete()

# Underline "||||" (this naturally separates it from "side" before it and
# "by_side|| Here is more text" after it, since the formatting differs).
$underlineRange = $d.Range($base + 24, $base + 28)
$underlineRange.Font.Underline = 1
